$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (shifts existing rows 8-22 down to 9-23)
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the data from the diff
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44741
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = 100112042
$ws.Cells.Item(8, 7).Value = "Locoto"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 2500
$ws.Cells.Item(8, 12).Value = 2500
$ws.Cells.Item(8, 13).Value = 2500
$ws.Cells.Item(8, 14).Value = "$/kilo"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 2500
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Copy style from row 9 (below it) for the D column date style
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 4).Value = 44741
